# Fix Paie tion Bug
# Update the district/municipality beneficiary counts and amounts in the two summary tables.
$d = $word.ActiveDocument

$cell = $d.Tables.Item(1).Cell(2,3)
$cell.Range.Text = "360"

$cell = $d.Tables.Item(1).Cell(2,4)
$cell.Range.Text = "3 600 000,00"

$cell = $d.Tables.Item(1).Cell(2,5)
$cell.Range.Text = "4 290 000,00"

$cell = $d.Tables.Item(1).Cell(3,3)
$cell.Range.Text = "69"

$cell = $d.Tables.Item(1).Cell(3,4)
$cell.Range.Text = "690 000,00"

$cell = $d.Tables.Item(1).Cell(4,3)
$cell.Range.Text = "61"

$cell = $d.Tables.Item(1).Cell(4,4)
$cell.Range.Text = "610 000,00"

$cell = $d.Tables.Item(1).Cell(4,5)
$cell.Range.Text = "1 640 000,00"

$cell = $d.Tables.Item(1).Cell(5,3)
$cell.Range.Text = "51"

$cell = $d.Tables.Item(1).Cell(5,4)
$cell.Range.Text = "510 000,00"

$cell = $d.Tables.Item(1).Cell(6,3)
$cell.Range.Text = "35"

$cell = $d.Tables.Item(1).Cell(6,4)
$cell.Range.Text = "350 000,00"

$cell = $d.Tables.Item(1).Cell(7,3)
$cell.Range.Text = "17"

$cell = $d.Tables.Item(1).Cell(7,4)
$cell.Range.Text = "170 000,00"

$cell = $d.Tables.Item(1).Cell(8,3)
$cell.Range.Text = "64"

$cell = $d.Tables.Item(1).Cell(8,4)
$cell.Range.Text = "640 000,00"

$cell = $d.Tables.Item(1).Cell(8,5)
$cell.Range.Text = "1 650 000,00"

$cell = $d.Tables.Item(1).Cell(9,3)
$cell.Range.Text = "44"

$cell = $d.Tables.Item(1).Cell(9,4)
$cell.Range.Text = "440 000,00"

$cell = $d.Tables.Item(1).Cell(10,3)
$cell.Range.Text = "23"

$cell = $d.Tables.Item(1).Cell(10,4)
$cell.Range.Text = "230 000,00"

$cell = $d.Tables.Item(1).Cell(11,3)
$cell.Range.Text = "16"

$cell = $d.Tables.Item(1).Cell(11,4)
$cell.Range.Text = "160 000,00"

$cell = $d.Tables.Item(1).Cell(12,3)
$cell.Range.Text = "18"

$cell = $d.Tables.Item(1).Cell(12,4)
$cell.Range.Text = "180 000,00"

$cell = $d.Tables.Item(1).Cell(13,3)
$cell.Range.Text = "134"

$cell = $d.Tables.Item(1).Cell(13,4)
$cell.Range.Text = "1 340 000,00"

$cell = $d.Tables.Item(1).Cell(13,5)
$cell.Range.Text = "1 700 000,00"

$cell = $d.Tables.Item(1).Cell(14,3)
$cell.Range.Text = "13"

$cell = $d.Tables.Item(1).Cell(14,4)
$cell.Range.Text = "130 000,00"

$cell = $d.Tables.Item(1).Cell(15,3)
$cell.Range.Text = "14"

$cell = $d.Tables.Item(1).Cell(15,4)
$cell.Range.Text = "140 000,00"

$cell = $d.Tables.Item(1).Cell(16,3)
$cell.Range.Text = "9"

$cell = $d.Tables.Item(1).Cell(16,4)
$cell.Range.Text = "90 000,00"

$cell = $d.Tables.Item(1).Cell(17,3)
$cell.Range.Text = "68"

$cell = $d.Tables.Item(1).Cell(17,4)
$cell.Range.Text = "680 000,00"

$cell = $d.Tables.Item(1).Cell(17,5)
$cell.Range.Text = "1 410 000,00"

$cell = $d.Tables.Item(1).Cell(18,3)
$cell.Range.Text = "8"

$cell = $d.Tables.Item(1).Cell(18,4)
$cell.Range.Text = "80 000,00"

$cell = $d.Tables.Item(1).Cell(19,3)
$cell.Range.Text = "37"

$cell = $d.Tables.Item(1).Cell(19,4)
$cell.Range.Text = "370 000,00"

$cell = $d.Tables.Item(1).Cell(20,3)
$cell.Range.Text = "28"

$cell = $d.Tables.Item(1).Cell(20,4)
$cell.Range.Text = "280 000,00"

$cell = $d.Tables.Item(1).Cell(21,3)
$cell.Range.Text = "27"

$cell = $d.Tables.Item(1).Cell(21,4)
$cell.Range.Text = "270 000,00"

$cell = $d.Tables.Item(1).Cell(21,5)
$cell.Range.Text = "1 100 000,00"

$cell = $d.Tables.Item(1).Cell(22,3)
$cell.Range.Text = "44"

$cell = $d.Tables.Item(1).Cell(22,4)
$cell.Range.Text = "440 000,00"

$cell = $d.Tables.Item(1).Cell(23,3)
$cell.Range.Text = "23"

$cell = $d.Tables.Item(1).Cell(23,4)
$cell.Range.Text = "230 000,00"

$cell = $d.Tables.Item(1).Cell(24,3)
$cell.Range.Text = "16"

$cell = $d.Tables.Item(1).Cell(24,4)
$cell.Range.Text = "160 000,00"

$cell = $d.Tables.Item(1).Cell(25,3)
$cell.Range.Text = "1179"

$cell = $d.Tables.Item(1).Cell(25,4)
$cell.Range.Text = "11 790 000,00"

$cell = $d.Tables.Item(1).Cell(25,5)
$cell.Range.Text = "11 790 000,00"

$cell = $d.Tables.Item(2).Cell(2,3)
$cell.Range.Text = "162"

$cell = $d.Tables.Item(2).Cell(2,4)
$cell.Range.Text = "1 620 000,00"

$cell = $d.Tables.Item(2).Cell(2,5)
$cell.Range.Text = "2 070 000,00"

$cell = $d.Tables.Item(2).Cell(3,3)
$cell.Range.Text = "25"

$cell = $d.Tables.Item(2).Cell(3,4)
$cell.Range.Text = "250 000,00"

$cell = $d.Tables.Item(2).Cell(4,3)
$cell.Range.Text = "20"

$cell = $d.Tables.Item(2).Cell(4,4)
$cell.Range.Text = "200 000,00"

$cell = $d.Tables.Item(2).Cell(5,3)
$cell.Range.Text = "95"

$cell = $d.Tables.Item(2).Cell(5,4)
$cell.Range.Text = "950 000,00"

$cell = $d.Tables.Item(2).Cell(5,5)
$cell.Range.Text = "1 130 000,00"

$cell = $d.Tables.Item(2).Cell(6,3)
$cell.Range.Text = "18"

$cell = $d.Tables.Item(2).Cell(6,4)
$cell.Range.Text = "180 000,00"

$cell = $d.Tables.Item(2).Cell(7,3)
$cell.Range.Text = "320"

$cell = $d.Tables.Item(2).Cell(7,4)
$cell.Range.Text = "3 200 000,00"

$cell = $d.Tables.Item(2).Cell(7,5)
$cell.Range.Text = "3 200 000,00"

$cell = $d.Tables.Item(2).Cell(8,3)
$cell.Range.Text = "1499"

$cell = $d.Tables.Item(2).Cell(8,4)
$cell.Range.Text = "14 990 000,00"

$cell = $d.Tables.Item(2).Cell(8,5)
$cell.Range.Text = "14 990 000,00"

# Update the amount-in-words paragraph to match the new grand total.
$d.Content.Find.Execute("واحد وأربعون مليون وتسعمئة ألف", $true, $false, $false, $false, $false, $true, 1, $false, "أربعة عشر مليون وتسعمئة وتسعون ألف", 2)
